# -----------------------------------------------------------------------
# edit.ps1 -- applies the "Statistical model and changed conclusion" diff
#
#   1. "4.2. Statistical Model" section: replaces the "To come." placeholder
#      with the full paragraph describing the logistic-regression model
#      (with an italicised "Sparkasse" in the middle of the text).
#   2. "5. Conclusion" section: rewords the sentence describing the sample
#      / methodology used in the analysis.
#
# Both edits are performed by locating the target paragraph with
# Range.Find, deleting its contents and re-inserting the replacement
# content as OOXML via Range.InsertXML (a WordOpenXML package, just like
# what Word puts on the clipboard). This lets us specify the exact run
# boundaries / run formatting required by the target document instead of
# relying on Word's automatic "merge adjacent runs with identical
# formatting" behaviour that a plain Range.Text assignment would trigger.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "4.2. Statistical Model" paragraph
# ---------------------------------------------------------------------------

$rng1 = $d.Content.Duplicate
$rng1.Find.Execute("To come.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng1.Find.Found) {
    $rng1.Delete()
    $xml1 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Given the research question, setting up the dependent variable as binary offers the best theoretical fit between the the construct under study and the data available to us. Therefore the analysis will employ a logistic regression model using re-election / no re-election of incumbent mayors as the dependent variable. This will allow for a very close fit of the model and the data, intuitive visualisations of patterns in the data using predictive probabilities and also keeps the complexity of the model in check. The primary indepent variable of the model is board membership in the local</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Sparkasse</w:t></w:r><w:r><w:t xml:space="preserve">, with additional controls for the top positions of chairman and vice-chairman. In order to cleanly estimate the effect of board membership on re-election chances, it is important to get a good grasp on the determinants of appointment to board positions. We therefore control for financial expertise of mayors, party affiliation as some party might contribute a disproportionate share of mayors. Moreover, we will control for other factors which might affect the chances of re-eleciton. These include the fiscal situation of the municipality, as municipalities in debt might reflect poorly on the incumbent, and the size of the municipality, as larger municipalities might offer more resources to retain the mayor''s office. For the size of the municipality, the number of valid votes is used as a proxy.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng1.InsertXML($xml1)
}

# ---------------------------------------------------------------------------
# Change 2: "5. Conclusion" paragraph
# ---------------------------------------------------------------------------

$oldPara2 = 'The paper purports to estimate the effect of board membership in German public savings banks (Sparkassen) on mayors'' electoral success. It compares electoral performance of mayors with and without a board seat in five German federal states covering 2,099 of 11,192 municipalities and 79 of 416 savings banks in Germany over the years from 2006 to 2015. To answer the research question, we compile a novel hand-collected dataset on the boardroom composition of German public banks with detailed information on board member profiles.'
$rng2 = $d.Content.Duplicate
$rng2.Find.Execute($oldPara2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
    $rng2.Delete()
    $xml2 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">The paper purports to estimate the effect of board membership in German public savings banks (</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Sparkassen</w:t></w:r><w:r><w:t xml:space="preserve">) on mayors'' electoral success. It compares electoral performance of mayors with and without a board seat in in Bavaria using a logistic regression model. To answer the research question, we compile a novel hand-collected dataset on the boardroom composition of German public banks with detailed information on board member profiles in addition to the existing data set on Bavarian municipal elections from 1948 to 2014.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng2.InsertXML($xml2)
}
